# edit.ps1 -- apply "Add files via upload" edit to the dev-log docx
#
# Summary of the target change (see unified diff):
#   1. Fix the typo "上船已經寫好的HTML網頁:" -> "上傳已經寫好的HTML網頁:"
#      (word was mistyped; the edit re-typed the first two characters,
#      which also relocates Word's hidden "_GoBack" last-edit bookmark
#      to right after "上傳").
#   2. Add a trailing full-width period to
#      "用拖拉或任意方法將相關檔案上傳" -> "用拖拉或任意方法將相關檔案上傳。"
#      and delete the now-redundant one-cell/one-row demo table that
#      held the literal text "index.html".
#
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "上船已經寫好的HTML網頁:" -> "上傳已經寫好的HTML網頁:"
#    Doing this as a Range.Text assignment (rather than Find/Replace)
#    reproduces Word's real behaviour of splitting the run at the edit
#    point and leaving the freshly-typed "_GoBack" bookmark there.
# ---------------------------------------------------------------------
$rngTitle = $d.Content
$rngTitle.Find.Execute("上船已經寫好的HTML網頁:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngTitle.Text = "上傳已經寫好的HTML網頁:"

$rngBookmark = $d.Content
$rngBookmark.Find.Execute("上傳", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBookmark.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngBookmark)

# ---------------------------------------------------------------------
# 2) Add the trailing "。" and remove the "index.html" demo table.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("用拖拉或任意方法將相關檔案上傳", $true, $false, $false, $false, $false, $true, 1, $false, "用拖拉或任意方法將相關檔案上傳。", 2)

$indexTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables($i)
    if ($candidate.Range.Text -like "*index.html*") {
        $indexTable = $candidate
        break
    }
}
if ($indexTable -ne $null) {
    $indexTable.Delete()
}
